# Add a new "step" row to the end of the RC-REF schema table.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Append a brand-new row after the last existing row (errorDistributionID).
$newRow = $t.Rows.Add()
$i = $newRow.Index

$t.Cell($i, 1).Range.Text = "step"
$t.Cell($i, 2).Range.Text = "Etape d'intégration du message"
$t.Cell($i, 3).Range.Text = "string"
$t.Cell($i, 4).Range.Text = "0..1"
$t.Cell($i, 5).Range.Text = "Nomenclature permettant d'identifier les différentes étapes d'intégration et de consultation du dossier dans le système émetteur"
# 6th cell ("Exemple") is left empty, matching the other rows in the table.
